$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "2025-05-27 08:58:42"
$ws.Range("B2").Value = 2536
$ws.Range("E2").Value = 1151
$ws.Range("F2").Value = 45.38643533123029
$ws.Range("G2").Value = 1719
$ws.Range("H2").Value = 67.78391167192429
$ws.Range("I2").Value = 804
$ws.Range("J2").Value = 31.70347003154574
$ws.Range("L2").Value = 13
$ws.Range("M2").Value = 0.5126182965299685
$ws.Range("T2").Value = 690
$ws.Range("U2").Value = 27.20820189274448
$ws.Range("W2").Value = 461
$ws.Range("X2").Value = 18.1782334384858
$ws.Range("Y2").Value = 114
$ws.Range("Z2").Value = 4.495268138801261
$ws.Range("AB2").Value = 1271
$ws.Range("AC2").Value = 50.11829652996845
$ws.Range("AD2").Value = 2536
$ws.Range("AE2").Value = 2424
$ws.Range("AF2").Value = 95.58359621451105
$ws.Range("AG2").Value = 112
$ws.Range("AH2").Value = 4.416403785488953
$ws.Range("AI2").Value = 126
$ws.Range("AJ2").Value = 174
$ws.Range("AK2").Value = 353
$ws.Range("AL2").Value = 19.29555895865237
$ws.Range("AM2").Value = 26.64624808575804
$ws.Range("AN2").Value = 54.05819295558959
$ws.Range("AO2").Value = 451674.65
$ws.Range("AP2").Value = 85459.66
$ws.Range("AQ2").Value = 28369.44
$ws.Range("AR2").Value = 79.87120332977456
$ws.Range("AS2").Value = 15.11212967199599
$ws.Range("AT2").Value = 5.016666998229454
$ws.Range("AU2").Value = 83.37217771303715
$ws.Range("AV2").Value = 130.1105769230769
$ws.Range("AW2").Value = 218.3128491620112
